$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plain text / non-numeric-looking values: assign directly
$ws.Range("D2").Value = "65.448.46"
$ws.Range("E2").Value = "  -1.04%  "
$ws.Range("D3").Value = "3.436.92"
$ws.Range("E3").Value = "  -3.22%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("E5").Value = "  -2.05%  "
$ws.Range("E6").Value = "  -7.34%  "
$ws.Range("D7").Value = "3.437.87"
$ws.Range("E7").Value = "  -3.14%  "
$ws.Range("E8").Value = "  -0.12%  "
$ws.Range("E9").Value = "  +0.46%  "
$ws.Range("E10").Value = "  -6.64%  "
$ws.Range("E11").Value = "  -8.80%  "
$ws.Range("E12").Value = "  -7.59%  "
$ws.Range("D13").Value = "4.017.29"
$ws.Range("E13").Value = "  -3.34%  "
$ws.Range("E14").Value = "  -11.00%  "
$ws.Range("E15").Value = "  -8.92%  "
$ws.Range("D16").Value = "3.444.49"
$ws.Range("E16").Value = "  -2.86%  "
$ws.Range("D17").Value = "65.467.28"
$ws.Range("E17").Value = "  -1.06%  "
$ws.Range("E18").Value = "  -2.13%  "
$ws.Range("E19").Value = "  -9.34%  "
$ws.Range("E20").Value = "  -6.65%  "
$ws.Range("E21").Value = "  -7.16%  "
$ws.Range("E22").Value = "  -5.90%  "
$ws.Range("E23").Value = "  -9.56%  "
$ws.Range("E24").Value = "  -5.66%  "
$ws.Range("E25").Value = "  +0.00%  "
$ws.Range("D26").Value = "3.579.14"
$ws.Range("E26").Value = "  -3.23%  "
$ws.Range("E27").Value = "  -10.48%  "
$ws.Range("E28").Value = "  -0.01%  "
$ws.Range("E29").Value = "  -8.67%  "
$ws.Range("E30").Value = "  -8.74%  "
$ws.Range("E31").Value = "  -11.24%  "
$ws.Range("D32").Value = "3.442.49"
$ws.Range("E32").Value = "  -3.03%  "
$ws.Range("E34").Value = "  -6.11%  "
$ws.Range("E35").Value = "  -6.72%  "
$ws.Range("E36").Value = "  -2.16%  "
$ws.Range("B37").Value = "Aptos"
$ws.Range("C37").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("E37").Value = "  -9.53%  "
$ws.Range("B38").Value = "Fetch.AI"
$ws.Range("C38").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("E38").Value = "  -12.78%  "
$ws.Range("E39").Value = "  -6.82%  "
$ws.Range("E40").Value = "  -10.77%  "
$ws.Range("E41").Value = "  -7.82%  "
$ws.Range("E42").Value = "  -5.41%  "
$ws.Range("E43").Value = "  -4.97%  "
$ws.Range("E44").Value = "  +0.04%  "
$ws.Range("E45").Value = "  -13.84%  "
$ws.Range("E46").Value = "  -11.58%  "
$ws.Range("E47").Value = "  -1.38%  "
$ws.Range("E48").Value = "  -3.20%  "
$ws.Range("E49").Value = "  -7.46%  "
$ws.Range("E50").Value = "  -15.48%  "
$ws.Range("D51").Value = "2.201.40"
$ws.Range("E51").Value = "  -7.29%  "

# Numeric-looking price values must stay as text: force text format,
# assign, then restore default style so no stray numeric formatting remains
$numericTextCells = @("D5", "D6", "D9", "D10", "D14", "D15", "D19", "D20", "D21", "D22", "D24", "D31", "D35", "D36", "D37", "D38", "D41", "D42", "D43", "D49")
foreach ($addr in $numericTextCells) {
    $ws.Range($addr).NumberFormat = "@"
}
$ws.Range("D5").Value = "592.20"
$ws.Range("D6").Value = "135.78"
$ws.Range("D9").Value = "0.489"
$ws.Range("D10").Value = "7.35"
$ws.Range("D14").Value = "0.0000180"
$ws.Range("D15").Value = "26.62"
$ws.Range("D19").Value = "9.92"
$ws.Range("D20").Value = "5.83"
$ws.Range("D21").Value = "13.64"
$ws.Range("D22").Value = "394.37"
$ws.Range("D24").Value = "73.29"
$ws.Range("D31").Value = "8.17"
$ws.Range("D35").Value = "23.00"
$ws.Range("D36").Value = "171.09"
$ws.Range("D37").Value = "6.93"
$ws.Range("D38").Value = "1.20"
$ws.Range("D41").Value = "0.0769"
$ws.Range("D42").Value = "0.823"
$ws.Range("D43").Value = "43.50"
$ws.Range("D49").Value = "6.55"
foreach ($addr in $numericTextCells) {
    $ws.Range($addr).Style = "Normal"
}
